$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.860.11'
$ws.Range("E2").Value = '  +4.01%  '

$ws.Range("D3").Value = '2.276.26'
$ws.Range("E3").Value = '  +4.42%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'251.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '

$ws.Range("D6").Value = "'0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.30%  '

$ws.Range("D7").Value = "'71.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.57%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").Value = "'0.640"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +11.16%  '

$ws.Range("D10").Value = "'38.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.40%  '

$ws.Range("D11").Value = "'59.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.66%  '

$ws.Range("D12").Value = "'0.0972"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.13%  '

$ws.Range("D13").Value = "'7.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.01%  '

$ws.Range("E14").Value = '  +1.76%  '

$ws.Range("D15").Value = '2.616.94'
$ws.Range("E15").Value = '  +4.27%  '

$ws.Range("D16").Value = "'14.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.93%  '

$ws.Range("D17").Value = "'0.887"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.81%  '

$ws.Range("D18").Value = '2.271.08'
$ws.Range("E18").Value = '  +4.50%  '

$ws.Range("D19").Value = '42.754.92'

$ws.Range("E20").Value = '  +7.45%  '

$ws.Range("E21").Value = '  +3.60%  '

$ws.Range("D22").Value = "'73.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.42%  '

$ws.Range("D23").Value = "'234.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.63%  '

$ws.Range("E24").Value = '  +4.94%  '

$ws.Range("D25").Value = "'4.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.19%  '

$ws.Range("D26").Value = "'11.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.54%  '

$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("E28").Value = '  +1.11%  '

$ws.Range("E29").Value = '  -1.25%  '

$ws.Range("D30").Value = "'2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.30%  '

$ws.Range("D31").Value = "'168.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").Value = "'21.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.70%  '

$ws.Range("D33").Value = "'6.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +11.91%  '

$ws.Range("D34").Value = "'0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.94%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.0799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.39%  '

$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = "'31.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +28.04%  '

$ws.Range("E37").Value = '  +3.76%  '

$ws.Range("D38").Value = "'4.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.91%  '

$ws.Range("E39").Value = '  +5.08%  '

$ws.Range("D40").Value = "'0.0312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.29%  '

$ws.Range("D41").Value = "'2.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.99%  '

$ws.Range("D42").Value = "'13.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +16.68%  '

$ws.Range("E43").Value = '  +5.88%  '

$ws.Range("D44").Value = "'0.211"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.86%  '

$ws.Range("D45").Value = "'5.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.03%  '

$ws.Range("E46").Value = '  +7.76%  '

$ws.Range("D47").Value = "'61.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("D48").Value = "'0.103"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.87%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = "'1.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.93%  '

$ws.Range("B50").Value = 'BinanceUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.01%  '

$ws.Range("E51").Value = '  +4.07%  '
